# Apply NATMI re-run values (3 ligand/receptor-expressing cells instead of 1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 101.814756
$ws.Range("H2").Value = 305.444268
$ws.Range("I2").Value = 0.2429953264021572
$ws.Range("J2").Value = 0.2429953264021571
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.193104333333333
$ws.Range("N2").Value = 3.579313
$ws.Range("O2").Value = 0.03883297235786565
$ws.Range("P2").Value = 0.03883297235786565
$ws.Range("Q2").Value = 121.475626580876
$ws.Range("R2").Value = 1093.280639227884
$ws.Range("S2").Value = 0.009436230793265511
$ws.Range("T2").Value = 0.00943623079326551

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 101.814756
$ws.Range("H3").Value = 305.444268
$ws.Range("I3").Value = 0.2429953264021572
$ws.Range("J3").Value = 0.2429953264021571
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 29.530898
$ws.Range("N3").Value = 88.59269400000001
$ws.Range("O3").Value = 0.9611670276421344
$ws.Range("P3").Value = 0.9611670276421344
$ws.Range("Q3").Value = 3006.681174330888
$ws.Range("R3").Value = 27060.13056897799
$ws.Range("S3").Value = 0.2335590956088917
$ws.Range("T3").Value = 0.2335590956088916

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 123.930687
$ws.Range("H4").Value = 371.792061
$ws.Range("I4").Value = 0.2957781260980341
$ws.Range("J4").Value = 0.2957781260980341
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.193104333333333
$ws.Range("N4").Value = 3.579313
$ws.Range("O4").Value = 0.03883297235786565
$ws.Range("P4").Value = 0.03883297235786565
$ws.Range("Q4").Value = 147.862239692677
$ws.Range("R4").Value = 1330.760157234093
$ws.Range("S4").Value = 0.01148594379482626
$ws.Range("T4").Value = 0.01148594379482626

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 123.930687
$ws.Range("H5").Value = 371.792061
$ws.Range("I5").Value = 0.2957781260980341
$ws.Range("J5").Value = 0.2957781260980341
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 29.530898
$ws.Range("N5").Value = 88.59269400000001
$ws.Range("O5").Value = 0.9611670276421344
$ws.Range("P5").Value = 0.9611670276421344
$ws.Range("Q5").Value = 3659.784476866926
$ws.Range("R5").Value = 32938.06029180234
$ws.Range("S5").Value = 0.2842921823032079
$ws.Range("T5").Value = 0.2842921823032079

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 115.753432
$ws.Range("H6").Value = 347.260296
$ws.Range("I6").Value = 0.2762619496039445
$ws.Range("J6").Value = 0.2762619496039445
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.193104333333333
$ws.Range("N6").Value = 3.579313
$ws.Range("O6").Value = 0.03883297235786565
$ws.Range("P6").Value = 0.03883297235786565
$ws.Range("Q6").Value = 138.1059213174053
$ws.Range("R6").Value = 1242.953291856648
$ws.Range("S6").Value = 0.01072807265250005
$ws.Range("T6").Value = 0.01072807265250005

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 115.753432
$ws.Range("H7").Value = 347.260296
$ws.Range("I7").Value = 0.2762619496039445
$ws.Range("J7").Value = 0.2762619496039445
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 29.530898
$ws.Range("N7").Value = 88.59269400000001
$ws.Range("O7").Value = 0.9611670276421344
$ws.Range("P7").Value = 0.9611670276421344
$ws.Range("Q7").Value = 3418.302793541936
$ws.Range("R7").Value = 30764.72514187742
$ws.Range("S7").Value = 0.2655338769514444
$ws.Range("T7").Value = 0.2655338769514444

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 77.499949
$ws.Range("H8").Value = 232.499847
$ws.Range("I8").Value = 0.1849645978958643
$ws.Range("J8").Value = 0.1849645978958643
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.193104333333333
$ws.Range("N8").Value = 3.579313
$ws.Range("O8").Value = 0.03883297235786565
$ws.Range("P8").Value = 0.03883297235786565
$ws.Range("Q8").Value = 92.46552498501234
$ws.Range("R8").Value = 832.1897248651111
$ws.Range("S8").Value = 0.007182725117273833
$ws.Range("T8").Value = 0.007182725117273833

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 77.499949
$ws.Range("H9").Value = 232.499847
$ws.Range("I9").Value = 0.1849645978958643
$ws.Range("J9").Value = 0.1849645978958643
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 29.530898
$ws.Range("N9").Value = 88.59269400000001
$ws.Range("O9").Value = 0.9611670276421344
$ws.Range("P9").Value = 0.9611670276421344
$ws.Range("Q9").Value = 2288.643088924202
$ws.Range("R9").Value = 20597.78780031782
$ws.Range("S9").Value = 0.1777818727785904
$ws.Range("T9").Value = 0.1777818727785904
